$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date serial numbers (Excel 1900 date system) replacing the inline-string
# "yyyy-mm-dd" text that used to live in column B, one entry per row from
# B2 to B296 (295 rows total).
$dateSerials = @(
    44348,44348,44351,44351,44351,44351,44351,44352,44352,44352,44353,44353,44353,44354,44354,44356,44356,44356,44357,44357,
    44358,44358,44358,44359,44359,44359,44360,44360,44360,44360,44360,44360,44361,44361,44362,44362,44362,44365,44365,44365,
    44365,44365,44366,44366,44366,44367,44367,44367,44367,44367,44367,44369,44369,44369,44370,44370,44370,44372,44372,44372,
    44372,44372,44373,44373,44373,44373,44373,44373,44374,44374,44374,44374,44374,44374,44375,44375,44376,44376,44376,44377,
    44377,44377,44378,44378,44379,44379,44379,44379,44380,44380,44380,44380,44380,44380,44381,44381,44381,44382,44382,44382,
    44383,44383,44383,44383,44384,44384,44384,44385,44385,44385,44386,44386,44386,44386,44386,44387,44387,44387,44387,44387,
    44387,44388,44388,44388,44388,44388,44389,44389,44389,44390,44390,44391,44391,44391,44392,44392,44392,44393,44393,44393,
    44394,44394,44394,44394,44394,44394,44395,44395,44395,44395,44395,44395,44396,44396,44397,44397,44397,44398,44398,44399,
    44399,44399,44400,44400,44400,44400,44400,44400,44401,44401,44401,44401,44401,44402,44402,44402,44402,44402,44403,44403,
    44404,44404,44405,44405,44405,44406,44406,44406,44407,44407,44407,44408,44408,44408,44408,44408,44408,44409,44409,44409,
    44409,44409,44409,44413,44413,44413,44414,44414,44414,44414,44414,44414,44415,44415,44415,44415,44415,44416,44416,44416,
    44416,44416,44416,44418,44418,44420,44420,44420,44421,44421,44421,44421,44421,44421,44422,44422,44422,44422,44423,44423,
    44423,44423,44423,44423,44424,44424,44424,44425,44425,44426,44426,44427,44427,44428,44428,44428,44428,44429,44429,44429,
    44429,44430,44430,44430,44430,44431,44431,44431,44432,44432,44432,44433,44433,44434,44434,44434,44435,44435,44435,44436,
    44436,44436,44437,44437,44437,44437,44437,44438,44438,44439,44439,44439,44440,44440,44440
)

# Mint custom number format 164 ("yyyy-mm-dd h:mm:ss") on B2 first, then
# immediately switch it to the uppercase variant, which becomes numFmtId
# 165 and the single cellXf actually used by every date cell in the column.
$ws.Cells.Item(2, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(2, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"

for ($i = 0; $i -lt $dateSerials.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 2)
    if ($row -gt 2) {
        $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    }
    $cell.Value = $dateSerials[$i]
}
